$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-09 Tuesday", "2024-07-10 Wednesday"),
    @("282×3=846", "588×9=5292"),
    @("296×7=2072", "337×9=3033"),
    @("569×7=3983", "816×5=4080"),
    @("867×9=7803", "926×3=2778"),
    @("378×2=756", "324×3=972"),
    @("478×2=956", "898×7=6286"),
    @("923×2=1846", "926×5=4630"),
    @("825×3=2475", "823×4=3292"),
    @("710×2=1420", "468×5=2340"),
    @("145×8=1160", "658×9=5922"),
    @("525×8=4200", "531×9=4779"),
    @("954×9=8586", "811×6=4866"),
    @("710×9=6390", "205×9=1845"),
    @("905×6=5430", "831×8=6648"),
    @("434×6=2604", "605×3=1815"),
    @("224×2=448", "986×9=8874"),
    @("936×4=3744", "910×4=3640"),
    @("784×6=4704", "778×9=7002"),
    @("668×7=4676", "417×2=834"),
    @("399×2=798", "325×7=2275"),
    @("900×4=3600", "598×3=1794"),
    @("640×2=1280", "987×9=8883"),
    @("356×9=3204", "818×5=4090"),
    @("250×5=1250", "540×9=4860"),
    @("563×2=1126", "319×3=957")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
